# Applies the "Added Nacubo GL Account Category" commit:
# Two new TestSteps rows ("checkAccessibility" / EditProfile.../ SavePopUp...)
# are inserted into the TestSteps sheet (sheet1), pushing the former rows
# 8-12 down to rows 9-11,13-14, and the active sheet/selection changes from
# TestData!G4 to TestSteps!C17.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)     # "TestSteps"

$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# STEP 1: copy cell formatting into place *before* touching any values,
# so every copy source still carries its original style.
# ---------------------------------------------------------------------

# Rows 10 and 13 need the "Normal 2" style (s=4) that column A of the
# original row 9 (jsClick) used.
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

$ws.Range("A9").Copy()
$ws.Range("A13").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

# Rows 9, 11 and 14 need the plain "Normal" style (s=2) that column A of
# row 2 (and most other data rows) uses.
$ws.Range("A2").Copy()
$ws.Range("A9").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

$ws.Range("A2").Copy()
$ws.Range("A11").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

$ws.Range("A2").Copy()
$ws.Range("A14").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

# Columns B and C of the brand-new rows 13/14 just need the common
# "Normal" data style (s=2).
$ws.Range("B2:C2").Copy()
$ws.Range("B13:C13").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

$ws.Range("B2:C2").Copy()
$ws.Range("B14:C14").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# STEP 2: write the final values for rows 8-14.
# ---------------------------------------------------------------------

# Row 8 (new step)
$ws.Range("A8").Value = "checkAccessibility"
$ws.Range("B8").Value = "EditProfile_AppropriationNumberProfile"
$ws.Range("C8").Value = $null

# Row 9 (was row 8)
$ws.Range("A9").Value = "enter_text"
$ws.Range("B9").Value = "txt_title"
$ws.Range("C9").Value = "getData=NewTitle"

# Row 10 (was row 9)
$ws.Range("A10").Value = "jsClick"
$ws.Range("B10").Value = "button_saveDetails"
$ws.Range("C10").Value = "getData=WaitForPageLoad"

# Row 11 (was row 10)
$ws.Range("A11").Value = "smallWaitForElementPresent"
$ws.Range("B11").Value = "button_save"
$ws.Range("C11").Value = "getData=WaitForPageLoad"

# Row 12 (new step)
$ws.Range("A12").Value = "checkAccessibility"
$ws.Range("B12").Value = "SavePopUp_AppropriationNumberProfile"
$ws.Range("C12").Value = $null

# Row 13 (was row 11)
$ws.Range("A13").Value = "jsClick"
$ws.Range("B13").Value = "button_save"
$ws.Range("C13").Value = "getData=WaitForPageLoad"

# Row 14 (was row 12)
$ws.Range("A14").Value = "smallWaitForElementPresent"
$ws.Range("B14").Value = "link_approNumberProfile_wait"
$ws.Range("C14").Value = "getData=WaitForPageLoad"

# ---------------------------------------------------------------------
# STEP 3: update the active sheet / selected cell to match the saved view
# (TestSteps tab active, C17 selected; TestData tab no longer active).
# ---------------------------------------------------------------------
$ws.Activate() | Out-Null
$ws.Range("C17").Select() | Out-Null
